# Update cryptocurrency price/symbol data on Sheet1 to reflect the
# latest scrape (GitHub Actions symbol-list update).
# Values in column D (Price) and the rotated Kick/BKEX/CEJI rows (B/C/D/E,
# rows 41-43) are stored as literal text in the workbook (inline strings),
# so we force each target cell to Text format before assigning the new
# value and then restore its style to Normal/General so no visible
# formatting change is introduced - this stops Excel's COM layer from
# "helpfully" re-interpreting numeric-looking text as a real number.

function Set-CellText($Ws, $CellRef, $NewText) {
    $cell = $Ws.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $NewText
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @('D2', '247.43'),
    @('D3', '22.41'),
    @('D4', '5.237'),
    @('D5', '0.05692'),
    @('D6', '3.418'),
    @('D7', '6.308'),
    @('D8', '0.8069'),
    @('D10', '0.1411'),
    @('D11', '0.07436'),
    @('D12', '0.03036'),
    @('D13', '0.03077'),
    @('D14', '0.09391'),
    @('D15', '3.890'),
    @('D16', '0.001573'),
    @('D17', '0.04796'),
    @('D18', '0.0005841'),
    @('D19', '0.006385'),
    @('D21', '0.0009959'),
    @('D22', '0.0001500'),
    @('D24', '2.192'),
    @('D27', '0.0004751'),
    @('D40', '0.03947'),
    @('B41', 'BKEXToken'),
    @('C41', 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'),
    @('D41', '0.1066'),
    @('E41', '40BKEXTokenBKK'),
    @('B42', 'CEJI'),
    @('C42', 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'),
    @('D42', '0.002681'),
    @('E42', '41CEJICEJI'),
    @('B43', 'KickToken'),
    @('C43', 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'),
    @('D43', '0.006831'),
    @('E43', '42KickTokenKICK'),
    @('D44', '0.008437'),
    @('D45', '0.00005582'),
    @('D47', '0.4501'),
    @('D48', '0.2021'),
    @('D50', '0.01010')
)

foreach ($u in $updates) {
    Set-CellText $ws $u[0] $u[1]
}

Write-Output ("Updated " + $updates.Count + " cells.")
